$p = $ppt.ActivePresentation
$nm = $p.NotesMaster
Write-Host "NotesMaster Height:" $nm.Height
Write-Host "NotesMaster Width:" $nm.Width

$hm = $p.HandoutMaster
Write-Host "HandoutMaster:" $hm

# try background of notes master
try {
  Write-Host "NM Background:" $nm.Background
} catch {
  Write-Host "NM Background failed:" $_.Exception.Message
}
